# Update odds values in row 2 of the active worksheet to match the
# FlashScore data refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value  = 2.5    # Odd_H_HT: 2.4 -> 2.5
$ws.Range("U2").Value  = 2.2    # Odd_BTTS_Yes: 2.1 -> 2.2
$ws.Range("V2").Value  = 1.62   # Odd_BTTS_No: 1.67 -> 1.62
$ws.Range("X2").Value  = 7      # Odd_CS_2-0: 7.5 -> 7
$ws.Range("AC2").Value = 7.5    # Odd_CS_0-0: 8 -> 7.5
$ws.Range("AF2").Value = 81     # Odd_CS_3-3: 67 -> 81
